# "Divide preparation into three stages"
#
# The video_name values are renamed from "<date>_<id>" to "<id>_<date>"
# (splitting date-prep from id-prep), the previously-missing target_frame /
# frame_length values for video 7 are filled in, the now-redundant
# num_of_frames column (F) is dropped, and frame_length (E) is given an
# explicit integer ("#") number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rename video_name from "<date>_<id>" to "<id>_<date>" ----------------
$ws.Range("B2").Value2 = "1_20180806"
$ws.Range("B3").Value2 = "2_20180920"
$ws.Range("B4").Value2 = "3_20180920"
$ws.Range("B5").Value2 = "4_20180920"
$ws.Range("B6").Value2 = "5_20180920"
$ws.Range("B7").Value2 = "6_20180920"
$ws.Range("B8").Value2 = "7_20180920"

# --- fill in the previously missing target_frame / frame_length for video 7
$ws.Range("D8").Value2 = "200, 215, 270, 465, 480, 590, 605"
$ws.Range("E8").Value2 = 15

# --- drop the num_of_frames column (F) entirely ----------------------------
$ws.Range("F1:F8").ClearContents()

# --- give frame_length (E) an explicit integer number format --------------
$ws.Range("E1:E8").NumberFormat = "#"

# --- column widths (best effort; Excel quantizes to whole pixels) ---------
$ws.Columns.Item(1).ColumnWidth = 7.666666666666667
$ws.Columns.Item(2).ColumnWidth = 18.833333333333332
$ws.Columns.Item(3).ColumnWidth = 7.666666666666667
$ws.Columns.Item(4).ColumnWidth = 56.333333333333336
$ws.Columns.Item(5).ColumnWidth = 7.666666666666667

# --- scroll the sheet back so column A is visible (was topLeftCell "C1") --
$excel.ActiveWindow.ScrollColumn = 1

# --- cosmetic workbook view tweak (tab ratio) ------------------------------
$wb.Windows.Item(1).TabRatio = 0.993

Write-Output "done"
